$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for each record.
# Update every data row (2-91) from serial date 45206 (2023-10-07)
# to serial date 45208 (2023-10-09).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45206) {
        $cell.Value = 45208
    }
}
